$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.071110984775541
$ws.Cells.Item(2, 4).Value = 1.080671682254996
$ws.Cells.Item(2, 5).Value = 1.065238116930853
$ws.Cells.Item(2, 6).Value = 1.087893628158011
$ws.Cells.Item(2, 9).Value = 1.047078001236949
$ws.Cells.Item(2, 10).Value = 1.076037234508205
$ws.Cells.Item(2, 11).Value = 1.083344506346445
$ws.Cells.Item(2, 12).Value = 1.067951896095224
$ws.Cells.Item(2, 13).Value = 1.090547718600062
$ws.Cells.Item(2, 14).Value = 1.077565330804402

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.07310119780953
$ws.Cells.Item(3, 4).Value = 1.082642305671393
$ws.Cells.Item(3, 5).Value = 1.066985779900424
$ws.Cells.Item(3, 6).Value = 1.090017158474467
$ws.Cells.Item(3, 9).Value = 1.047596280181536
$ws.Cells.Item(3, 10).Value = 1.077680372855335
$ws.Cells.Item(3, 11).Value = 1.085130639702144
$ws.Cells.Item(3, 12).Value = 1.069512589696057
$ws.Cells.Item(3, 13).Value = 1.092487780279063
$ws.Cells.Item(3, 14).Value = 1.079210802596456

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.074383986401248
$ws.Cells.Item(4, 4).Value = 1.083912730216477
$ws.Cells.Item(4, 5).Value = 1.068111773981607
$ws.Cells.Item(4, 6).Value = 1.091386654773028
$ws.Cells.Item(4, 9).Value = 1.047927900307911
$ws.Cells.Item(4, 10).Value = 1.078738353684974
$ws.Cells.Item(4, 11).Value = 1.086281275986368
$ws.Cells.Item(4, 12).Value = 1.070517154658589
$ws.Cells.Item(4, 13).Value = 1.093738175423016
$ws.Cells.Item(4, 14).Value = 1.080270285880226

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.074922097258765
$ws.Cells.Item(5, 4).Value = 1.084445716840047
$ws.Cells.Item(5, 5).Value = 1.068584001522928
$ws.Cells.Item(5, 6).Value = 1.091961324227641
$ws.Cells.Item(5, 9).Value = 1.048066426137761
$ws.Cells.Item(5, 10).Value = 1.079181896962761
$ws.Cells.Item(5, 11).Value = 1.086763803752896
$ws.Cells.Item(5, 12).Value = 1.070938224961158
$ws.Cells.Item(5, 13).Value = 1.094262681001888
$ws.Cells.Item(5, 14).Value = 1.080714459040341

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.075012380392874
$ws.Cells.Item(6, 4).Value = 1.084535143879804
$ws.Cells.Item(6, 5).Value = 1.068663224406758
$ws.Cells.Item(6, 6).Value = 1.092057751988088
$ws.Cells.Item(6, 9).Value = 1.048089633444816
$ws.Cells.Item(6, 10).Value = 1.079256298244075
$ws.Cells.Item(6, 11).Value = 1.086844752628189
$ws.Cells.Item(6, 12).Value = 1.071008851922467
$ws.Cells.Item(6, 13).Value = 1.094350680481319
$ws.Cells.Item(6, 14).Value = 1.080788965980009

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.074391181240067
$ws.Cells.Item(7, 4).Value = 1.083919856297148
$ws.Cells.Item(7, 5).Value = 1.068118088358572
$ws.Cells.Item(7, 6).Value = 1.09139433768944
$ws.Cells.Item(7, 9).Value = 1.047929754771916
$ws.Cells.Item(7, 10).Value = 1.078744285141987
$ws.Cells.Item(7, 11).Value = 1.086287728227788
$ws.Cells.Item(7, 12).Value = 1.070522785897896
$ws.Cells.Item(7, 13).Value = 1.093745188417867
$ws.Cells.Item(7, 14).Value = 1.080276225760588

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.071784640467511
$ws.Cells.Item(8, 4).Value = 1.081338652376608
$ws.Cells.Item(8, 5).Value = 1.06582976836583
$ws.Cells.Item(8, 6).Value = 1.088612247872047
$ws.Cells.Item(8, 9).Value = 1.047253935684353
$ws.Cells.Item(8, 10).Value = 1.076593639315816
$ws.Cells.Item(8, 11).Value = 1.083949211343426
$ws.Cells.Item(8, 12).Value = 1.068480452201953
$ws.Cells.Item(8, 13).Value = 1.091204412876374
$ws.Cells.Item(8, 14).Value = 1.078122525770661

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.067151995260362
$ws.Cells.Item(9, 4).Value = 1.076753087707624
$ws.Cells.Item(9, 5).Value = 1.061759174754538
$ws.Cells.Item(9, 6).Value = 1.083673618590441
$ws.Cells.Item(9, 9).Value = 1.046034026429075
$ws.Cells.Item(9, 10).Value = 1.072762779450461
$ws.Cells.Item(9, 11).Value = 1.079788213855876
$ws.Cells.Item(9, 12).Value = 1.064839950698121
$ws.Cells.Item(9, 13).Value = 1.08668813882116
$ws.Cells.Item(9, 14).Value = 1.074286225644892

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.064035259753549
$ws.Cells.Item(10, 4).Value = 1.073669440421155
$ws.Cells.Item(10, 5).Value = 1.059018222536608
$ws.Cells.Item(10, 6).Value = 1.080355103131849
$ws.Cells.Item(10, 9).Value = 1.045200701411179
$ws.Cells.Item(10, 10).Value = 1.070179772437328
$ws.Cells.Item(10, 11).Value = 1.076985637535481
$ws.Cells.Item(10, 12).Value = 1.062383546520295
$ws.Cells.Item(10, 13).Value = 1.08364936803436
$ws.Cells.Item(10, 14).Value = 1.071699550465524

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.062678585718862
$ws.Cells.Item(11, 4).Value = 1.072327511974273
$ws.Cells.Item(11, 5).Value = 1.057824569412894
$ws.Cells.Item(11, 6).Value = 1.078911569729347
$ws.Cells.Item(11, 9).Value = 1.044834985820811
$ws.Cells.Item(11, 10).Value = 1.069054076943547
$ws.Cells.Item(11, 11).Value = 1.075764969220191
$ws.Cells.Item(11, 12).Value = 1.061312607559926
$ws.Cells.Item(11, 13).Value = 1.082326557686778
$ws.Cells.Item(11, 14).Value = 1.070572256355024

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.062173555903316
$ws.Cells.Item(12, 4).Value = 1.071828023407887
$ws.Cells.Item(12, 5).Value = 1.057380142512671
$ws.Cells.Item(12, 6).Value = 1.07837435260263
$ws.Cells.Item(12, 9).Value = 1.044698399024002
$ws.Cells.Item(12, 10).Value = 1.068634828298993
$ws.Cells.Item(12, 11).Value = 1.075310456990384
$ws.Cells.Item(12, 12).Value = 1.060913689112946
$ws.Cells.Item(12, 13).Value = 1.081834123169161
$ws.Cells.Item(12, 14).Value = 1.070152412329306

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.062281936855087
$ws.Cells.Item(13, 4).Value = 1.07193521282272
$ws.Cells.Item(13, 5).Value = 1.057475521620201
$ws.Cells.Item(13, 6).Value = 1.078489634411654
$ws.Cells.Item(13, 9).Value = 1.044727731182376
$ws.Cells.Item(13, 10).Value = 1.068724809477051
$ws.Cells.Item(13, 11).Value = 1.075408001743168
$ws.Cells.Item(13, 12).Value = 1.060999309782201
$ws.Cells.Item(13, 13).Value = 1.08193980160754
$ws.Cells.Item(13, 14).Value = 1.070242521290958

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.06263686244255
$ws.Cells.Item(14, 4).Value = 1.072286245424768
$ws.Cells.Item(14, 5).Value = 1.057787854555535
$ws.Cells.Item(14, 6).Value = 1.078867184292165
$ws.Cells.Item(14, 9).Value = 1.044823710747249
$ws.Cells.Item(14, 10).Value = 1.069019444633654
$ws.Cells.Item(14, 11).Value = 1.075727421745425
$ws.Cells.Item(14, 12).Value = 1.061279655926264
$ws.Cells.Item(14, 13).Value = 1.082285875192097
$ws.Cells.Item(14, 14).Value = 1.070537574863278

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.062855396966986
$ws.Cells.Item(15, 4).Value = 1.072502389881715
$ws.Cells.Item(15, 5).Value = 1.057980153045923
$ws.Cells.Item(15, 6).Value = 1.079099668450478
$ws.Cells.Item(15, 9).Value = 1.044882748053749
$ws.Cells.Item(15, 10).Value = 1.069200830454539
$ws.Cells.Item(15, 11).Value = 1.075924080068949
$ws.Cells.Item(15, 12).Value = 1.061452236667795
$ws.Cells.Item(15, 13).Value = 1.082498957967713
$ws.Cells.Item(15, 14).Value = 1.070719218272834

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.064125145206342
$ws.Cells.Item(16, 4).Value = 1.07375835617497
$ws.Cells.Item(16, 5).Value = 1.059097295576955
$ws.Cells.Item(16, 6).Value = 1.080450763879989
$ws.Cells.Item(16, 9).Value = 1.045224869032705
$ws.Cells.Item(16, 10).Value = 1.070254326294937
$ws.Cells.Item(16, 11).Value = 1.077066496425736
$ws.Cells.Item(16, 12).Value = 1.062454465091943
$ws.Cells.Item(16, 13).Value = 1.083737008361974
$ws.Cells.Item(16, 14).Value = 1.071774210198163

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.064919698905674
$ws.Cells.Item(17, 4).Value = 1.074544377660121
$ws.Cells.Item(17, 5).Value = 1.05979620799476
$ws.Cells.Item(17, 6).Value = 1.081296481151438
$ws.Cells.Item(17, 9).Value = 1.045438158500299
$ws.Cells.Item(17, 10).Value = 1.070913199364686
$ws.Cells.Item(17, 11).Value = 1.077781173072116
$ws.Cells.Item(17, 12).Value = 1.063081162997276
$ws.Cells.Item(17, 13).Value = 1.084511707592875
$ws.Cells.Item(17, 14).Value = 1.072434018943245

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.065382465163932
$ws.Cells.Item(18, 4).Value = 1.075002207483809
$ws.Cells.Item(18, 5).Value = 1.060203217396933
$ws.Cells.Item(18, 6).Value = 1.081789139783535
$ws.Cells.Item(18, 9).Value = 1.045562096474425
$ws.Cells.Item(18, 10).Value = 1.071296812585671
$ws.Cells.Item(18, 11).Value = 1.078197345552608
$ws.Cells.Item(18, 12).Value = 1.063446002830638
$ws.Cells.Item(18, 13).Value = 1.084962903076655
$ws.Cells.Item(18, 14).Value = 1.072818176939004

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.065540141559102
$ws.Cells.Item(19, 4).Value = 1.075158207473688
$ws.Cells.Item(19, 5).Value = 1.060341886959066
$ws.Cells.Item(19, 6).Value = 1.081957017133505
$ws.Cells.Item(19, 9).Value = 1.045604276688563
$ws.Cells.Item(19, 10).Value = 1.071427497518625
$ws.Cells.Item(19, 11).Value = 1.078339134097936
$ws.Cells.Item(19, 12).Value = 1.063570285458163
$ws.Cells.Item(19, 13).Value = 1.085116635799301
$ws.Cells.Item(19, 14).Value = 1.072949047459553

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.06483452174488
$ws.Cells.Item(20, 4).Value = 1.074460111751689
$ws.Cells.Item(20, 5).Value = 1.059721289201317
$ws.Cells.Item(20, 6).Value = 1.081205809507146
$ws.Cells.Item(20, 9).Value = 1.045415323265912
$ws.Cells.Item(20, 10).Value = 1.070842580742321
$ws.Cells.Item(20, 11).Value = 1.077704566235013
$ws.Cells.Item(20, 12).Value = 1.063013997097562
$ws.Cells.Item(20, 13).Value = 1.084428659586104
$ws.Cells.Item(20, 14).Value = 1.072363300034337

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.062532376354338
$ws.Cells.Item(21, 4).Value = 1.07218290394459
$ws.Cells.Item(21, 5).Value = 1.057695909537795
$ws.Cells.Item(21, 6).Value = 1.078756033819399
$ws.Cells.Item(21, 9).Value = 1.044795467773983
$ws.Cells.Item(21, 10).Value = 1.068932712936376
$ws.Cells.Item(21, 11).Value = 1.075633391169377
$ws.Cells.Item(21, 12).Value = 1.061197132224731
$ws.Cells.Item(21, 13).Value = 1.08218399532477
$ws.Cells.Item(21, 14).Value = 1.070450719997042

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.06107854132562
$ws.Cells.Item(22, 4).Value = 1.0707451195555
$ws.Cells.Item(22, 5).Value = 1.056416377620987
$ws.Cells.Item(22, 6).Value = 1.077209817746779
$ws.Cells.Item(22, 9).Value = 1.04440143161127
$ws.Cells.Item(22, 10).Value = 1.067725436732581
$ws.Cells.Item(22, 11).Value = 1.07432477217194
$ws.Cells.Item(22, 12).Value = 1.060048280003532
$ws.Cells.Item(22, 13).Value = 1.080766398829621
$ws.Cells.Item(22, 14).Value = 1.069241729322561

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.061849863122882
$ws.Cells.Item(23, 4).Value = 1.071507896976914
$ws.Cells.Item(23, 5).Value = 1.057095269251942
$ws.Cells.Item(23, 6).Value = 1.078030071081642
$ws.Cells.Item(23, 9).Value = 1.044610729585202
$ws.Cells.Item(23, 10).Value = 1.068366059116709
$ws.Cells.Item(23, 11).Value = 1.075019111550635
$ws.Cells.Item(23, 12).Value = 1.060657935292459
$ws.Cells.Item(23, 13).Value = 1.081518500238014
$ws.Cells.Item(23, 14).Value = 1.069883261463955

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.064873011748176
$ws.Cells.Item(24, 4).Value = 1.074498189874932
$ws.Cells.Item(24, 5).Value = 1.059755143800258
$ws.Cells.Item(24, 6).Value = 1.081246782077978
$ws.Cells.Item(24, 9).Value = 1.045425642977837
$ws.Cells.Item(24, 10).Value = 1.070874492408567
$ws.Cells.Item(24, 11).Value = 1.07773918368685
$ws.Cells.Item(24, 12).Value = 1.063044348646101
$ws.Cells.Item(24, 13).Value = 1.084466187485497
$ws.Cells.Item(24, 14).Value = 1.072395257018808

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.068354508725117
$ws.Cells.Item(25, 4).Value = 1.077943134690414
$ws.Cells.Item(25, 5).Value = 1.062816208833517
$ws.Cells.Item(25, 6).Value = 1.084954843379055
$ws.Cells.Item(25, 9).Value = 1.046352896832358
$ws.Cells.Item(25, 10).Value = 1.073758171130946
$ws.Cells.Item(25, 11).Value = 1.080868855054144
$ws.Cells.Item(25, 12).Value = 1.065786186122342
$ws.Cells.Item(25, 13).Value = 1.087860504200513
$ws.Cells.Item(25, 14).Value = 1.075283030895736
